# Refresh market-price derived columns (H-N) per scheduled runner update.
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1369.68
$ws.Range("J129").Value = 1441.3043
$ws.Range("L129").Value = 4323.9129
$ws.Range("N129").Value = -14323.9129
# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 223785
$ws.Range("I132").Value = 1580.1052
$ws.Range("J132").Value = 1430040.1
$ws.Range("K132").Value = 4740.3156
$ws.Range("L132").Value = 4290120.300000001
$ws.Range("M132").Value = -2210.3156
$ws.Range("N132").Value = -4295180.300000001
# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2704.12
$ws.Range("I138").Value = 691
$ws.Range("J138").Value = 3448.6987
$ws.Range("K138").Value = 2073
$ws.Range("L138").Value = 10346.0961
$ws.Range("M138").Value = 3067
$ws.Range("N138").Value = -20626.0961

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1278.84
$ws.Range("I61").Value = 1261.1
$ws.Range("J61").Value = 1349.8
$ws.Range("K61").Value = 1261.1
$ws.Range("L61").Value = 1349.8
$ws.Range("M61").Value = -1049.1
$ws.Range("N61").Value = -1773.8
# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3550
$ws.Range("N122").ClearContents()
# row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1278.84
$ws.Range("I136").Value = 1261.1
$ws.Range("J136").Value = 1349.8
$ws.Range("K136").Value = 3783.3
$ws.Range("L136").Value = 4049.4
$ws.Range("M136").Value = -1233.3
$ws.Range("N136").Value = -9149.4

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 853.5714
$ws.Range("I94").Value = 394.55554
$ws.Range("K94").Value = 394.55554
$ws.Range("M94").Value = 56.44445999999999
# row 137 (Leve Item ID 42153)
$ws.Range("H137").Value = 49117.65
$ws.Range("J137").Value = 49117.65
$ws.Range("L137").Value = 49117.65
$ws.Range("N137").Value = -59317.65

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 136.72728
$ws.Range("I7").Value = 141.44444
$ws.Range("J7").Value = 115.5
$ws.Range("K7").Value = 141.44444
$ws.Range("L7").Value = 115.5
$ws.Range("M7").Value = -28.44443999999999
$ws.Range("N7").Value = -341.5
# row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 1121.9474
$ws.Range("I16").Value = 1093.6
$ws.Range("J16").Value = 1228.25
$ws.Range("K16").Value = 1093.6
$ws.Range("L16").Value = 1228.25
$ws.Range("M16").Value = -806.5999999999999
$ws.Range("N16").Value = -1802.25
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1752.1632
$ws.Range("I31").Value = 1762.6364
$ws.Range("J31").Value = 1660
$ws.Range("K31").Value = 1762.6364
$ws.Range("L31").Value = 1660
$ws.Range("M31").Value = -1467.6364
$ws.Range("N31").Value = -2250
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1752.1632
$ws.Range("I34").Value = 1762.6364
$ws.Range("J34").Value = 1660
$ws.Range("K34").Value = 1762.6364
$ws.Range("L34").Value = 1660
$ws.Range("M34").Value = -1560.6364
$ws.Range("N34").Value = -2064
# row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 1121.9474
$ws.Range("I113").Value = 1093.6
$ws.Range("J113").Value = 1228.25
$ws.Range("K113").Value = 1093.6
$ws.Range("L113").Value = 1228.25
$ws.Range("M113").Value = 1076.4
$ws.Range("N113").Value = -5568.25
# row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1995.4324
$ws.Range("I134").Value = 1416.5938
$ws.Range("J134").Value = 5700
$ws.Range("K134").Value = 4249.7814
$ws.Range("L134").Value = 17100
$ws.Range("M134").Value = -1714.7814
$ws.Range("N134").Value = -22170

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row 8 (Leve Item ID 16734)
$ws.Range("H8").Value = 102.55556
$ws.Range("I8").Value = 102.55556
$ws.Range("K8").Value = 307.66668
$ws.Range("M8").Value = -168.66668
# row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 1345.9286
$ws.Range("I114").Value = 786.7143
$ws.Range("J114").Value = 1532.3334
$ws.Range("K114").Value = 2360.1429
$ws.Range("L114").Value = 4597.0002
$ws.Range("M114").Value = 893.8571000000002
$ws.Range("N114").Value = -11105.0002
# row 118 (Leve Item ID 27872)
$ws.Range("H118").Value = 4573.6875
$ws.Range("I118").Value = 912.375
$ws.Range("J118").Value = 8235
$ws.Range("K118").Value = 2737.125
$ws.Range("L118").Value = 24705
$ws.Range("M118").Value = -1494.125
$ws.Range("N118").Value = -27191

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row 3 (Leve Item ID 4091)
$ws.Range("H3").Value = 8750248
$ws.Range("I3").Value = 11666668
$ws.Range("J3").Value = 990
$ws.Range("K3").Value = 11666668
$ws.Range("L3").Value = 990
$ws.Range("M3").Value = -11666552
$ws.Range("N3").Value = -1222
# row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 4658.3335
$ws.Range("J70").Value = 6014.2856
$ws.Range("L70").Value = 6014.2856
$ws.Range("N70").Value = -6554.2856
# row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 4658.3335
$ws.Range("J73").Value = 6014.2856
$ws.Range("L73").Value = 6014.2856
$ws.Range("N73").Value = -7886.2856
# row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 4250
$ws.Range("I80").Value = 4250
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4250
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3252
$ws.Range("N80").ClearContents()
# row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 4250
$ws.Range("I83").Value = 4250
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 21250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -16258
$ws.Range("N83").ClearContents()
# row 119 (Leve Item ID 26282)
$ws.Range("H119").Value = 24750
$ws.Range("J119").Value = 24750
$ws.Range("L119").Value = 24750
$ws.Range("N119").Value = -34426

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 807.5
$ws.Range("I22").Value = 714.2857
$ws.Range("J22").Value = 938
$ws.Range("K22").Value = 714.2857
$ws.Range("L22").Value = 938
$ws.Range("M22").Value = -419.2857
$ws.Range("N22").Value = -1528
# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 807.5
$ws.Range("I27").Value = 714.2857
$ws.Range("J27").Value = 938
$ws.Range("K27").Value = 714.2857
$ws.Range("L27").Value = 938
$ws.Range("M27").Value = -607.2857
$ws.Range("N27").Value = -1152
# row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 1224.0555
$ws.Range("I61").Value = 944.8570999999999
$ws.Range("J61").Value = 2201.25
$ws.Range("K61").Value = 944.8570999999999
$ws.Range("L61").Value = 2201.25
$ws.Range("M61").Value = -742.8570999999999
$ws.Range("N61").Value = -2605.25
# row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 2079.2
$ws.Range("I93").Value = 896
$ws.Range("J93").Value = 2375
$ws.Range("K93").Value = 896
$ws.Range("L93").Value = 2375
$ws.Range("M93").Value = 352
$ws.Range("N93").Value = -4871
# row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 1224.0555
$ws.Range("I113").Value = 944.8570999999999
$ws.Range("J113").Value = 2201.25
$ws.Range("K113").Value = 944.8570999999999
$ws.Range("L113").Value = 2201.25
$ws.Range("M113").Value = 1225.1429
$ws.Range("N113").Value = -6541.25

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 1412.7142
$ws.Range("I122").Value = 1372.5
$ws.Range("J122").Value = 1466.3334
$ws.Range("K122").Value = 4117.5
$ws.Range("L122").Value = 4399.0002
$ws.Range("M122").Value = -1667.5
$ws.Range("N122").Value = -9299.0002
